$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 16563.875
$ws.Range("I20").Value = 16563.875
$ws.Range("K20").Value = 16563.875
$ws.Range("M20").Value = -16333.875
$ws.Range("H35").Value = 16563.875
$ws.Range("I35").Value = 16563.875
$ws.Range("K35").Value = 16563.875
$ws.Range("M35").Value = -16184.875
$ws.Range("H70").Value = 46770.59
$ws.Range("I70").Value = 501150
$ws.Range("J70").Value = 1332.65
$ws.Range("K70").Value = 1503450
$ws.Range("L70").Value = 3997.95
$ws.Range("M70").Value = -1503180
$ws.Range("N70").Value = -4537.950000000001
$ws.Range("H73").Value = 46770.59
$ws.Range("I73").Value = 501150
$ws.Range("J73").Value = 1332.65
$ws.Range("K73").Value = 1503450
$ws.Range("L73").Value = 3997.95
$ws.Range("M73").Value = -1502514
$ws.Range("N73").Value = -5869.950000000001
$ws.Range("H125").Value = 143390.14
$ws.Range("I125").Value = 167038.5
$ws.Range("J125").Value = 1500
$ws.Range("K125").Value = 1503346.5
$ws.Range("L125").Value = 13500
$ws.Range("M125").Value = -1500886.5
$ws.Range("N125").Value = -18420
$ws.Range("H132").Value = 2571.5
$ws.Range("I132").Value = 1961.7037
$ws.Range("J132").Value = 3540
$ws.Range("K132").Value = 5885.1111
$ws.Range("L132").Value = 10620
$ws.Range("M132").Value = -3355.1111
$ws.Range("N132").Value = -15680
$ws.Range("H137").Value = 1648.12
$ws.Range("I137").Value = 1094.3334
$ws.Range("J137").Value = 3072.1428
$ws.Range("K137").Value = 3283.0002
$ws.Range("L137").Value = 9216.428400000001
$ws.Range("M137").Value = -733.0001999999999
$ws.Range("N137").Value = -14316.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7917.8887
$ws.Range("I32").Value = 7738.7925
$ws.Range("K32").Value = 7738.7925
$ws.Range("M32").Value = -7451.7925
$ws.Range("H61").Value = 2765.0312
$ws.Range("I61").Value = 1458.68
$ws.Range("J61").Value = 7430.5713
$ws.Range("K61").Value = 1458.68
$ws.Range("L61").Value = 7430.5713
$ws.Range("M61").Value = -1246.68
$ws.Range("N61").Value = -7854.5713
$ws.Range("H123").Value = 24294.6
$ws.Range("J123").Value = 24294.6
$ws.Range("L123").Value = 24294.6
$ws.Range("N123").Value = -34094.6
$ws.Range("H136").Value = 2765.0312
$ws.Range("I136").Value = 1458.68
$ws.Range("J136").Value = 7430.5713
$ws.Range("K136").Value = 4376.04
$ws.Range("L136").Value = 22291.7139
$ws.Range("M136").Value = -1826.04
$ws.Range("N136").Value = -27391.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2359.4878
$ws.Range("I31").Value = 2263.5652
$ws.Range("J31").Value = 2482.0557
$ws.Range("K31").Value = 2263.5652
$ws.Range("L31").Value = 2482.0557
$ws.Range("M31").Value = -1968.5652
$ws.Range("N31").Value = -3072.0557
$ws.Range("H34").Value = 2359.4878
$ws.Range("I34").Value = 2263.5652
$ws.Range("J34").Value = 2482.0557
$ws.Range("K34").Value = 2263.5652
$ws.Range("L34").Value = 2482.0557
$ws.Range("M34").Value = -2061.5652
$ws.Range("N34").Value = -2886.0557
$ws.Range("H58").Value = 588936.6
$ws.Range("I58").Value = 823950.4
$ws.Range("J58").Value = 1402.1111
$ws.Range("K58").Value = 823950.4
$ws.Range("L58").Value = 1402.1111
$ws.Range("M58").Value = -823747.4
$ws.Range("N58").Value = -1808.1111
$ws.Range("H86").Value = 14107.7
$ws.Range("I86").Value = 36002.332
$ws.Range("J86").Value = 4724.2856
$ws.Range("K86").Value = 36002.332
$ws.Range("L86").Value = 4724.2856
$ws.Range("M86").Value = -34879.332
$ws.Range("N86").Value = -6970.2856
$ws.Range("H89").Value = 14107.7
$ws.Range("I89").Value = 36002.332
$ws.Range("J89").Value = 4724.2856
$ws.Range("K89").Value = 180011.66
$ws.Range("L89").Value = 23621.428
$ws.Range("M89").Value = -174395.66
$ws.Range("N89").Value = -34853.428
$ws.Range("H134").Value = 1119.3788
$ws.Range("I134").Value = 975.4888999999999
$ws.Range("J134").Value = 1427.7142
$ws.Range("K134").Value = 2926.4667
$ws.Range("L134").Value = 4283.142599999999
$ws.Range("M134").Value = -391.4666999999999
$ws.Range("N134").Value = -9353.142599999999
$ws.Range("H136").Value = 588936.6
$ws.Range("I136").Value = 823950.4
$ws.Range("J136").Value = 1402.1111
$ws.Range("K136").Value = 2471851.2
$ws.Range("L136").Value = 4206.3333
$ws.Range("M136").Value = -2469301.2
$ws.Range("N136").Value = -9306.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 33335266
$ws.Range("I22").Value = 50000400
$ws.Range("J22").Value = 4997.5
$ws.Range("K22").Value = 150001200
$ws.Range("L22").Value = 14992.5
$ws.Range("M22").Value = -150001031
$ws.Range("N22").Value = -15330.5
$ws.Range("H27").Value = 33335266
$ws.Range("I27").Value = 50000400
$ws.Range("J27").Value = 4997.5
$ws.Range("K27").Value = 150001200
$ws.Range("L27").Value = 14992.5
$ws.Range("M27").Value = -150001098
$ws.Range("N27").Value = -15196.5
$ws.Range("H75").Value = 11570.714
$ws.Range("J75").Value = 11570.714
$ws.Range("L75").Value = 34712.142
$ws.Range("N75").Value = -36708.142
$ws.Range("H78").Value = 11570.714
$ws.Range("J78").Value = 11570.714
$ws.Range("L78").Value = 104136.426
$ws.Range("N78").Value = -114120.426
$ws.Range("H134").Value = 4245.0884
$ws.Range("I134").Value = 2356.923
$ws.Range("J134").Value = 5413.952
$ws.Range("K134").Value = 7070.768999999999
$ws.Range("L134").Value = 16241.856
$ws.Range("M134").Value = -2000.768999999999
$ws.Range("N134").Value = -26381.856
$ws.Range("H140").Value = 2065
$ws.Range("I140").Value = 977
$ws.Range("K140").Value = 2931
$ws.Range("M140").Value = 2249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 30500
$ws.Range("J62").Value = 30500
$ws.Range("L62").Value = 30500
$ws.Range("N62").Value = -31872
$ws.Range("H65").Value = 30500
$ws.Range("J65").Value = 30500
$ws.Range("L65").Value = 91500
$ws.Range("N65").Value = -98364
$ws.Range("H109").Value = 9277.23
$ws.Range("J109").Value = 9277.23
$ws.Range("L109").Value = 9277.23
$ws.Range("N109").Value = -11357.23
$ws.Range("H123").Value = 8625
$ws.Range("J123").Value = 8625
$ws.Range("L123").Value = 8625
$ws.Range("N123").Value = -13525

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 4754
$ws.Range("I9").Value = 420.2857
$ws.Range("J9").Value = 14866
$ws.Range("K9").Value = 420.2857
$ws.Range("L9").Value = 14866
$ws.Range("M9").Value = -196.2857
$ws.Range("N9").Value = -15314
$ws.Range("H13").Value = 13000
$ws.Range("I13").Value = 13000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 13000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -12860
$ws.Range("N13").ClearContents()
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 45000
$ws.Range("L141").Value = 45000
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1225.5807
$ws.Range("I136").Value = 1246.4333
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 3739.2999
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = -1189.2999
$ws.Range("N136").Value = -6900
